# Add a new "number_of_run" parameter column to the "scenarios" sheet.
# This inserts a new column D (shifting the existing D:I columns to E:J),
# fills in the header + values for the new column, and restores the
# "no fill" formatting on the new header cell so it matches the rest of
# the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("scenarios")

# Shift columns D:I -> E:J, leaving a blank column D with C's formatting.
$ws.Columns("D:D").Insert()

# New column header + values (agent run count per scenario).
$ws.Range("D1").Value = "number_of_run"
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 2
$ws.Range("D4").Value = 1

# The inserted header cell picked up column C's fill; clear it back to
# "no fill" like the rest of row 1.
$ws.Range("D1").Interior.ColorIndex = -4142

# Column D was sized to fit the new "number_of_run" header text (stored
# column width of 14 once Excel applies its character/pixel padding).
$ws.Columns("D:D").ColumnWidth = 13.3

# Move the active selection, as recorded in the saved workbook.
$ws.Range("E5").Select() | Out-Null
